# Auto-generated Excel COM-interop script to apply value updates
# described by the Excalibur_Profits.xlsx diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(8, 8).Value = 247.22223
$ws.Cells.Item(8, 9).Value = 174.71428
$ws.Cells.Item(8, 10).Value = 501
$ws.Cells.Item(8, 11).Value = 524.14284
$ws.Cells.Item(8, 12).Value = 1503
$ws.Cells.Item(8, 13).Value = -385.14284
$ws.Cells.Item(8, 14).Value = -1781

$ws.Cells.Item(17, 8).Value = 1782.8823
$ws.Cells.Item(17, 10).Value = 1833.75
$ws.Cells.Item(17, 12).Value = 5501.25
$ws.Cells.Item(17, 14).Value = -5837.25

$ws.Cells.Item(18, 8).Value = 1310.2778
$ws.Cells.Item(18, 9).Value = 891.1539
$ws.Cells.Item(18, 10).Value = 2400
$ws.Cells.Item(18, 11).Value = 891.1539
$ws.Cells.Item(18, 12).Value = 2400
$ws.Cells.Item(18, 13).Value = -607.1539
$ws.Cells.Item(18, 14).Value = -2968

$ws.Cells.Item(33, 8).Value = 738.375
$ws.Cells.Item(33, 9).Value = 393.5
$ws.Cells.Item(33, 10).Value = 1773
$ws.Cells.Item(33, 11).Value = 393.5
$ws.Cells.Item(33, 12).Value = 1773
$ws.Cells.Item(33, 13).Value = -164.5
$ws.Cells.Item(33, 14).Value = -2231

$ws.Cells.Item(93, 8).Value = 89799.57000000001
$ws.Cells.Item(93, 10).Value = 89799.57000000001
$ws.Cells.Item(93, 12).Value = 89799.57000000001
$ws.Cells.Item(93, 14).Value = -94791.57000000001

$ws.Cells.Item(99, 8).Value = 1178.8182
$ws.Cells.Item(99, 9).Value = 209.2
$ws.Cells.Item(99, 11).Value = 627.5999999999999
$ws.Cells.Item(99, 13).Value = 870.4000000000001

$ws.Cells.Item(125, 8).Value = 5851.7144
$ws.Cells.Item(125, 9).Value = 5997.4
$ws.Cells.Item(125, 10).Value = 5770.778
$ws.Cells.Item(125, 11).Value = 53976.6
$ws.Cells.Item(125, 12).Value = 51937.002
$ws.Cells.Item(125, 13).Value = -51516.6
$ws.Cells.Item(125, 14).Value = -56857.002

$ws.Cells.Item(138, 8).Value = 4632.4287
$ws.Cells.Item(138, 10).Value = 5691.6924
$ws.Cells.Item(138, 12).Value = 17075.0772
$ws.Cells.Item(138, 14).Value = -27355.0772

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 628.5341
$ws.Cells.Item(32, 9).Value = 566.7931
$ws.Cells.Item(32, 11).Value = 566.7931
$ws.Cells.Item(32, 13).Value = -279.7931

$ws.Cells.Item(45, 8).Value = 3206.4285
$ws.Cells.Item(45, 9).Value = 3315.3333
$ws.Cells.Item(45, 11).Value = 3315.3333
$ws.Cells.Item(45, 13).Value = -2938.3333

$ws.Cells.Item(102, 8).Value = 64739
$ws.Cells.Item(102, 9).Value = 74712.336
$ws.Cells.Item(102, 11).Value = 74712.336
$ws.Cells.Item(102, 13).Value = -73090.336

$ws.Cells.Item(123, 8).Value = 45000
$ws.Cells.Item(123, 10).Value = 45000
$ws.Cells.Item(123, 12).Value = 45000
$ws.Cells.Item(123, 14).Value = -54800

$ws.Cells.Item(132, 8).Value = 2513.6316
$ws.Cells.Item(132, 9).Value = 2523.2646
$ws.Cells.Item(132, 10).Value = 2431.75
$ws.Cells.Item(132, 11).Value = 7569.793799999999
$ws.Cells.Item(132, 12).Value = 7295.25
$ws.Cells.Item(132, 13).Value = -5039.793799999999
$ws.Cells.Item(132, 14).Value = -12355.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(9, 8).Value = 50000
$ws.Cells.Item(9, 10).Value = 50000
$ws.Cells.Item(9, 12).Value = 50000
$ws.Cells.Item(9, 14).Value = -50336

$ws.Cells.Item(99, 8).Value = 2945.8462
$ws.Cells.Item(99, 9).Value = 2572.3635
$ws.Cells.Item(99, 11).Value = 2572.3635
$ws.Cells.Item(99, 13).Value = -1074.3635

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 98.18519000000001
$ws.Cells.Item(7, 9).Value = 74.36842
$ws.Cells.Item(7, 10).Value = 154.75
$ws.Cells.Item(7, 11).Value = 74.36842
$ws.Cells.Item(7, 12).Value = 154.75
$ws.Cells.Item(7, 13).Value = 38.63158
$ws.Cells.Item(7, 14).Value = -380.75

$ws.Cells.Item(16, 8).Value = 1250
$ws.Cells.Item(16, 9).Value = 1250
$ws.Cells.Item(16, 11).Value = 1250
$ws.Cells.Item(16, 13).Value = -963

$ws.Cells.Item(31, 8).Value = 14860.788
$ws.Cells.Item(31, 9).Value = 1074.92
$ws.Cells.Item(31, 10).Value = 57941.625
$ws.Cells.Item(31, 11).Value = 1074.92
$ws.Cells.Item(31, 12).Value = 57941.625
$ws.Cells.Item(31, 13).Value = -779.9200000000001
$ws.Cells.Item(31, 14).Value = -58531.625

$ws.Cells.Item(34, 8).Value = 14860.788
$ws.Cells.Item(34, 9).Value = 1074.92
$ws.Cells.Item(34, 10).Value = 57941.625
$ws.Cells.Item(34, 11).Value = 1074.92
$ws.Cells.Item(34, 12).Value = 57941.625
$ws.Cells.Item(34, 13).Value = -872.9200000000001
$ws.Cells.Item(34, 14).Value = -58345.625

$ws.Cells.Item(113, 8).Value = 1250
$ws.Cells.Item(113, 9).Value = 1250
$ws.Cells.Item(113, 11).Value = 1250
$ws.Cells.Item(113, 13).Value = 920

$ws.Cells.Item(120, 8).Value = 50000
$ws.Cells.Item(120, 10).Value = 50000
$ws.Cells.Item(120, 12).Value = 50000
$ws.Cells.Item(120, 14).Value = -57258

$ws.Cells.Item(127, 8).Value = 100000
$ws.Cells.Item(127, 10).Value = 100000
$ws.Cells.Item(127, 12).Value = 100000
$ws.Cells.Item(127, 14).Value = -109920

$ws.Cells.Item(132, 8).Value = 2828.0688
$ws.Cells.Item(132, 9).Value = 2482
$ws.Cells.Item(132, 11).Value = 7446
$ws.Cells.Item(132, 13).Value = -4916

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 2511417
$ws.Cells.Item(7, 9).Value = 4018111
$ws.Cells.Item(7, 10).Value = 260
$ws.Cells.Item(7, 11).Value = 12054333
$ws.Cells.Item(7, 12).Value = 780
$ws.Cells.Item(7, 13).Value = -12054221
$ws.Cells.Item(7, 14).Value = -1004

$ws.Cells.Item(11, 8).Value = 2000327.6
$ws.Cells.Item(11, 9).Value = 350
$ws.Cells.Item(11, 10).Value = 2500322
$ws.Cells.Item(11, 11).Value = 1050
$ws.Cells.Item(11, 12).Value = 7500966
$ws.Cells.Item(11, 13).Value = -910
$ws.Cells.Item(11, 14).Value = -7501246

$ws.Cells.Item(12, 8).Value = 284.5
$ws.Cells.Item(12, 10).Value = 305.77274
$ws.Cells.Item(12, 12).Value = 917.31822
$ws.Cells.Item(12, 14).Value = -1263.31822

$ws.Cells.Item(97, 8).Value = 275
$ws.Cells.Item(97, 9).Value = 323.33334
$ws.Cells.Item(97, 10).Value = 250.83333
$ws.Cells.Item(97, 11).Value = 970.0000200000001
$ws.Cells.Item(97, 12).Value = 752.49999
$ws.Cells.Item(97, 13).Value = -474.0000200000001
$ws.Cells.Item(97, 14).Value = -1744.49999

$ws.Cells.Item(118, 8).Value = 654.5
$ws.Cells.Item(118, 9).Value = 385.4
$ws.Cells.Item(118, 10).Value = 2000
$ws.Cells.Item(118, 11).Value = 1156.2
$ws.Cells.Item(118, 12).Value = 6000
$ws.Cells.Item(118, 13).Value = 86.80000000000018
$ws.Cells.Item(118, 14).Value = -8486

$ws.Cells.Item(137, 8).Value = 2109.25
$ws.Cells.Item(137, 9).Value = 1026.3334
$ws.Cells.Item(137, 10).Value = 2759
$ws.Cells.Item(137, 11).Value = 3079.0002
$ws.Cells.Item(137, 12).Value = 8277
$ws.Cells.Item(137, 13).Value = 2020.9998
$ws.Cells.Item(137, 14).Value = -18477

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(52, 8).Value = 45000
$ws.Cells.Item(52, 10).Value = 45000
$ws.Cells.Item(52, 12).Value = 45000
$ws.Cells.Item(52, 14).Value = -45518

$ws.Cells.Item(62, 8).Value = 75762.86
$ws.Cells.Item(62, 9).Value = 56666.668
$ws.Cells.Item(62, 11).Value = 56666.668
$ws.Cells.Item(62, 13).Value = -55980.668

$ws.Cells.Item(65, 8).Value = 75762.86
$ws.Cells.Item(65, 9).Value = 56666.668
$ws.Cells.Item(65, 11).Value = 170000.004
$ws.Cells.Item(65, 13).Value = -166568.004

$ws.Cells.Item(93, 8).Value = 55403
$ws.Cells.Item(93, 10).Value = 55403
$ws.Cells.Item(93, 12).Value = 55403
$ws.Cells.Item(93, 14).Value = -59147

$ws.Cells.Item(97, 8).Value = 2276.95
$ws.Cells.Item(97, 9).Value = 1522
$ws.Cells.Item(97, 11).Value = 1522
$ws.Cells.Item(97, 13).Value = -1026

$ws.Cells.Item(107, 8).Value = 1384.9286
$ws.Cells.Item(107, 10).Value = 1732.5
$ws.Cells.Item(107, 12).Value = 1732.5
$ws.Cells.Item(107, 14).Value = -5572.5

$ws.Cells.Item(122, 8).Value = 3830.6191
$ws.Cells.Item(122, 9).Value = 7812
$ws.Cells.Item(122, 11).Value = 23436
$ws.Cells.Item(122, 13).Value = -20986

$ws.Cells.Item(132, 8).Value = 4478.069
$ws.Cells.Item(132, 9).Value = 4057.8572
$ws.Cells.Item(132, 11).Value = 12173.5716
$ws.Cells.Item(132, 13).Value = -9643.571599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(32, 8).Value = 4783.1
$ws.Cells.Item(32, 9).Value = 4783.1
$ws.Cells.Item(32, 11).Value = 4783.1
$ws.Cells.Item(32, 13).Value = -4466.1

$ws.Cells.Item(35, 8).Value = 1357.8125
$ws.Cells.Item(35, 9).Value = 1357.8125
$ws.Cells.Item(35, 11).Value = 1357.8125
$ws.Cells.Item(35, 13).Value = -1021.8125

$ws.Cells.Item(48, 8).Value = 0
$ws.Cells.Item(48, 9).Value = 0
$ws.Cells.Item(48, 11).Value = 0
$ws.Cells.Item(48, 13).ClearContents()

$ws.Cells.Item(55, 8).Value = 688.5
$ws.Cells.Item(55, 9).Value = 665.0769
$ws.Cells.Item(55, 10).Value = 749.4
$ws.Cells.Item(55, 11).Value = 665.0769
$ws.Cells.Item(55, 12).Value = 749.4
$ws.Cells.Item(55, 13).Value = -492.0769
$ws.Cells.Item(55, 14).Value = -1095.4

$ws.Cells.Item(93, 8).Value = 1504.8422
$ws.Cells.Item(93, 9).Value = 1424.5
$ws.Cells.Item(93, 10).Value = 1933.3334
$ws.Cells.Item(93, 11).Value = 1424.5
$ws.Cells.Item(93, 12).Value = 1933.3334
$ws.Cells.Item(93, 13).Value = -176.5
$ws.Cells.Item(93, 14).Value = -4429.3334

$ws.Cells.Item(100, 8).Value = 9800.429
$ws.Cells.Item(100, 9).Value = 2654.3635
$ws.Cells.Item(100, 11).Value = 2654.3635
$ws.Cells.Item(100, 13).Value = -2113.3635

$ws.Cells.Item(122, 8).Value = 36310.25
$ws.Cells.Item(122, 9).Value = 4763
$ws.Cells.Item(122, 11).Value = 14289
$ws.Cells.Item(122, 13).Value = -11839

$ws.Cells.Item(132, 8).Value = 29039
$ws.Cells.Item(132, 9).Value = 32159.334
$ws.Cells.Item(132, 10).Value = 14997.5
$ws.Cells.Item(132, 11).Value = 96478.00199999999
$ws.Cells.Item(132, 12).Value = 44992.5
$ws.Cells.Item(132, 13).Value = -93948.00199999999
$ws.Cells.Item(132, 14).Value = -50052.5

$ws.Cells.Item(136, 8).Value = 24310.62
$ws.Cells.Item(136, 9).Value = 1474.2354
$ws.Cells.Item(136, 11).Value = 4422.706200000001
$ws.Cells.Item(136, 13).Value = -1872.706200000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 3154.58
$ws.Cells.Item(132, 9).Value = 1975.2106
$ws.Cells.Item(132, 11).Value = 5925.6318
$ws.Cells.Item(132, 13).Value = -3395.6318

$ws.Cells.Item(136, 8).Value = 8692.272000000001
$ws.Cells.Item(136, 9).Value = 8177.775
$ws.Cells.Item(136, 11).Value = 24533.325
$ws.Cells.Item(136, 13).Value = -21983.325
